$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 1")

# Column E (index 5) gets a numeric suffix appended, cycling 1/2/3 as the
# category repeats across the table; column I (index 9) switches from
# "Aftermarket" to "Genuine" for every existing data row.
$suffixes = @{
    2  = "1"; 3  = "1"; 4  = "1"; 5  = "1"; 6  = "1"; 7  = "1"; 8  = "1"
    9  = "2"; 10 = "2"; 11 = "2"; 12 = "2"; 13 = "2"; 14 = "2"; 15 = "2"
    16 = "3"; 17 = "3"; 18 = "3"; 19 = "3"; 20 = "3"; 21 = "3"; 22 = "3"
}

for ($r = 2; $r -le 22; $r++) {
    $base = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 5).Value = [string]$base + $suffixes[$r]
}

for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 9).Value = "Genuine"
}

# Append a new row 23 duplicating the pattern of row 2 (Honda Sedan Sensor)
# but with D23 = "S2" and E23 left without a numeric suffix.
$ws.Cells.Item(23, 1).Formula = '= (10*10+ROW()-1)&"-S"&ROW()-1&"-K"&ROW()-1&"-U"&ROW()-1'
$ws.Cells.Item(23, 2).Value = "Honda"
$ws.Cells.Item(23, 3).Value = "Sedan"
$ws.Cells.Item(23, 4).Value = "S2"
$ws.Cells.Item(23, 5).Value = "Sensor"
$ws.Cells.Item(23, 6).Formula = '=(20*10+ROW()-1)&"-I"&ROW()-1&"-C"'
$ws.Cells.Item(23, 7).Formula = '=B23&" "&C23&" "&E23'
$ws.Cells.Item(23, 8).Value = "Research Oem"
$ws.Cells.Item(23, 9).Value = "Genuine"

$ws.Range("K21").Select()
